$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsDeDe = $wb.Worksheets.Item(3)

# The handoff failed, so the status text changes everywhere it is shown
# ("Ready for handoff" -> "Handoff transform failed").
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"
$wsZhCn.Range("B2").Value = "Handoff transform failed"
$wsDeDe.Range("B2").Value = "Handoff transform failed"

# --- zh-cn sheet: the handoff transform failed, so there is no handoff file,
#     the handoff timestamp is reset, and the reason becomes "Ignored" ---
$wsZhCn.Range("C2").Hyperlinks.Delete()
$wsZhCn.Range("C2").Clear()
$wsZhCn.Range("D2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("H2").Value = "Ignored"

# Re-create the two remaining hyperlinks that were on this sheet (deleting
# any hyperlink clears the whole worksheet's hyperlink collection).
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5b79ba9bd48395576784d13279ba71cf7caa9345/e2e/f60e3c32-a450-403e-a589-e5ae29063f69.md", "", "", "f60e3c32-a450-403e-a589-e5ae29063f69.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5b79ba9bd48395576784d13279ba71cf7caa9345/.localization-config", "", "", ".localization-config")

# --- de-de sheet: same transform-failed handling ---
$wsDeDe.Range("C2").Hyperlinks.Delete()
$wsDeDe.Range("C2").Clear()
$wsDeDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("H2").Value = "Ignored"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5b79ba9bd48395576784d13279ba71cf7caa9345/e2e/f60e3c32-a450-403e-a589-e5ae29063f69.md", "", "", "f60e3c32-a450-403e-a589-e5ae29063f69.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5b79ba9bd48395576784d13279ba71cf7caa9345/.localization-config", "", "", ".localization-config")

Write-Host "Handoff report generated"
